$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8; existing rows 8-58 shift down to 9-59.
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the new weekly record.
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = 44532
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 100114007
$ws.Range("G8").Value = "Jengibre"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 220
$ws.Range("K8").Value = 13000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 13909
$ws.Range("N8").Value = '$/caja 13 kilos'
$ws.Range("O8").Value = "Perú"
$ws.Range("P8").Value = 1070
$ws.Range("Q8").Value = 13
$ws.Range("R8").Value = "Hortaliza"
